$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("year" variable): label changes from "YEAR" to "year"
$ws.Range("E10").Value = "year"

# Row 28 ("pasture_nr" variable): update label to reflect NASS inflation-adjusted rent
$ws.Range("E28").Value = "2010USD pastureland rent/acre (NASS)"

# Row 31 ("CRPrent" variable): update label to reflect inflation-adjusted CRP payments
$ws.Range("E31").Value = "2010USD CRP Contract-based FY rent payments (not actuals) (USDA)"

# Row 32 ("CRP_nr" variable): update label to reflect inflation-adjusted avg CRP payments
$ws.Range("E32").Value = "2010USD avg per-CRPacre contract-based FY rent payments (not actuals) (USDA)"
